$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 5885
$ws.Range("J3").Value = 6288
$ws.Range("J4").Value = 1360
$ws.Range("J5").Value = 481
$ws.Range("J6").Value = 8079
$ws.Range("J7").Value = 22093

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 175
$ws.Range("J6").Value = 161
$ws.Range("J7").Value = 654
$ws.Range("J8").Value = 1388
$ws.Range("J9").Value = 110
$ws.Range("J10").Value = 154
$ws.Range("J11").Value = 349
$ws.Range("J14").Value = 108
$ws.Range("J15").Value = 244
$ws.Range("J16").Value = 87
$ws.Range("J18").Value = 183
$ws.Range("J19").Value = 656
$ws.Range("J20").Value = 456
$ws.Range("J23").Value = 208
$ws.Range("J27").Value = 134
$ws.Range("J29").Value = 1220
$ws.Range("J33").Value = 1026
$ws.Range("J36").Value = 306
$ws.Range("J37").Value = 678
$ws.Range("J41").Value = 143
$ws.Range("J42").Value = 929
$ws.Range("J43").Value = 178
$ws.Range("J44").Value = 169
$ws.Range("J47").Value = 168
$ws.Range("J48").Value = 259
$ws.Range("J50").Value = 132
$ws.Range("J51").Value = 272
$ws.Range("J52").Value = 554
$ws.Range("J53").Value = 306
$ws.Range("J54").Value = 432
$ws.Range("J55").Value = 305
$ws.Range("J63").Value = 78
$ws.Range("J64").Value = 145
$ws.Range("J65").Value = 556
$ws.Range("J67").Value = 834
$ws.Range("J73").Value = 212
$ws.Range("J75").Value = 68
$ws.Range("J76").Value = 336
$ws.Range("J77").Value = 165
$ws.Range("J78").Value = 272
$ws.Range("J79").Value = 631
$ws.Range("J83").Value = 447
$ws.Range("J85").Value = 910
$ws.Range("J86").Value = 141
$ws.Range("J87").Value = 74
$ws.Range("J89").Value = 289
$ws.Range("J90").Value = 237
$ws.Range("J94").Value = 226
$ws.Range("J96").Value = 253
$ws.Range("J97").Value = 190
$ws.Range("J101").Value = 22093

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("J6").Value = 38
$ws.Range("J7").Value = 108

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J6").Value = 89
$ws.Range("J7").Value = 253

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J2").Value = 200
$ws.Range("J6").Value = 213
$ws.Range("J7").Value = 654

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J2").Value = 105
$ws.Range("J6").Value = 146
$ws.Range("J7").Value = 349

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J3").Value = 80
$ws.Range("J7").Value = 289

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 239
$ws.Range("J4").Value = 60
$ws.Range("J6").Value = 267
$ws.Range("J7").Value = 910

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J2").Value = 132
$ws.Range("J6").Value = 224
$ws.Range("J7").Value = 554

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J6").Value = 201
$ws.Range("J7").Value = 306

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 380
$ws.Range("J3").Value = 424
$ws.Range("J6").Value = 471
$ws.Range("J7").Value = 1388

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J3").Value = 164
$ws.Range("J7").Value = 447

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J3").Value = 343
$ws.Range("J5").Value = 43
$ws.Range("J6").Value = 353
$ws.Range("J7").Value = 1026

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 202
$ws.Range("J7").Value = 678

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J5").Value = 15
$ws.Range("J6").Value = 197
$ws.Range("J7").Value = 556

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J2").Value = 208
$ws.Range("J5").Value = 25
$ws.Range("J7").Value = 834

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J3").Value = 86
$ws.Range("J6").Value = 209
$ws.Range("J7").Value = 432

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 372
$ws.Range("J3").Value = 426
$ws.Range("J7").Value = 1220

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J3").Value = 46
$ws.Range("J7").Value = 259

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J6").Value = 254
$ws.Range("J7").Value = 656

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J6").Value = 70
$ws.Range("J7").Value = 169

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J4").Value = 26
$ws.Range("J6").Value = 188
$ws.Range("J7").Value = 336

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("J4").Value = 11
$ws.Range("J7").Value = 161

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("J6").Value = 80
$ws.Range("J7").Value = 143

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J2").Value = 201
$ws.Range("J6").Value = 479
$ws.Range("J7").Value = 929

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J6").Value = 81
$ws.Range("J7").Value = 154

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J4").Value = 30
$ws.Range("J7").Value = 272

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J6").Value = 158
$ws.Range("J7").Value = 305

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J6").Value = 57
$ws.Range("J7").Value = 208

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 177
$ws.Range("J3").Value = 215
$ws.Range("J6").Value = 184
$ws.Range("J7").Value = 631

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J2").Value = 40
$ws.Range("J3").Value = 37
$ws.Range("J7").Value = 145

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J2").Value = 125
$ws.Range("J3").Value = 158
$ws.Range("J6").Value = 123
$ws.Range("J7").Value = 456

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("J2").Value = 49
$ws.Range("J7").Value = 183

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J6").Value = 94
$ws.Range("J7").Value = 306

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J2").Value = 41
$ws.Range("J6").Value = 124
$ws.Range("J7").Value = 226

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("J3").Value = 44
$ws.Range("J7").Value = 168

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J4").Value = 10
$ws.Range("J7").Value = 244

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("J6").Value = 41
$ws.Range("J7").Value = 132

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("J2").Value = 29
$ws.Range("J6").Value = 37
$ws.Range("J7").Value = 110

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J3").Value = 54
$ws.Range("J7").Value = 212

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J6").Value = 66
$ws.Range("J7").Value = 175

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J2").Value = 30
$ws.Range("J6").Value = 131
$ws.Range("J7").Value = 190

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J2").Value = 38
$ws.Range("J3").Value = 32
$ws.Range("J4").Value = 17
$ws.Range("J7").Value = 134

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J4").Value = 75
$ws.Range("J7").Value = 141

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("J3").Value = 21
$ws.Range("J7").Value = 68

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J2").Value = 84
$ws.Range("J5").Value = 8
$ws.Range("J6").Value = 70
$ws.Range("J7").Value = 237

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J3").Value = 74
$ws.Range("J7").Value = 272

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J6").Value = 103
$ws.Range("J7").Value = 178

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("J2").Value = 64
$ws.Range("J7").Value = 165

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("J6").Value = 49
$ws.Range("J7").Value = 74

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("J6").Value = 69
$ws.Range("J7").Value = 87
